# Update countries & provincias Spain
# Applies the 30-Jun-2020 04:33 -> 05:50 refresh of the "paises" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Refresh timestamp banner (A1)
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 05:50"

# -----------------------------------------------------------------
# 2) Re-order country rows whose case counts changed enough to move
#    them past their neighbour in the (descending) sort on column B,
#    by swapping the country-name text between the two rows.
# -----------------------------------------------------------------

# Honduras overtakes Japon (rows 55/56)
$nameJapon    = $ws.Range("A55").Text
$nameHonduras = $ws.Range("A56").Text
$ws.Range("A55").Value = $nameHonduras
$ws.Range("A56").Value = $nameJapon

# Butan overtakes San Martin (Parte Holandesa) (rows 186/187)
$nameSanMartin = $ws.Range("A186").Text
$nameButan     = $ws.Range("A187").Text
$ws.Range("A186").Value = $nameButan
$ws.Range("A187").Value = $nameSanMartin

# Dominica overtakes Fiyi (rows 205/206) - tied totals, order only
$nameFiyi     = $ws.Range("A205").Text
$nameDominica = $ws.Range("A206").Text
$ws.Range("A205").Value = $nameDominica
$ws.Range("A206").Value = $nameFiyi

# Groenlandia overtakes Islas Malvinas (rows 209/210) - tied totals, order only
$nameMalvinas   = $ws.Range("A209").Text
$nameGroenlandia = $ws.Range("A210").Text
$ws.Range("A209").Value = $nameGroenlandia
$ws.Range("A210").Value = $nameMalvinas

# -----------------------------------------------------------------
# 3) Updated per-country figures
#    Columns: B=Casos totales C=Nuevos casos D=Casos activos
#             E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes
# -----------------------------------------------------------------

# Row 55 - now Honduras
$ws.Range("B55").Value = 18818
$ws.Range("C55").Value = 736
$ws.Range("D55").Value = 1961
$ws.Range("E55").Value = 16372
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 6
$ws.Range("H55").Value = 485

# Row 56 - now Japon
$ws.Range("B56").Value = 18476
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 16557
$ws.Range("E56").Value = 947
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 972

# Row 75 - Australia
$ws.Range("B75").Value = 7836
$ws.Range("C75").Value = 69
$ws.Range("D75").Value = 7008
$ws.Range("E75").Value = 724
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 104

# Row 186 - now Butan
$ws.Range("B186").Value = 77
$ws.Range("C186").Value = 1
$ws.Range("D186").Value = 44
$ws.Range("E186").Value = 33
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

# Row 187 - now San Martin (Parte Holandesa)
$ws.Range("B187").Value = 77
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 62
$ws.Range("E187").Value = 0
$ws.Range("F187").Value = 0
$ws.Range("G187").Value = 0
$ws.Range("H187").Value = 15
